$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: condition object declarations (drop leading "$")
$ws.Range("B6").Value = "form:RequirementForm"
$ws.Range("C6").Value = "caseData:CaseData"

# Row 7: condition / action expressions (switch from getter-call syntax to bare field access)
$ws.Range("B7").Value = "code in (`$param)"
$ws.Range("C7").Value = "caseSource in (`$param)"
$ws.Range("D7").Value = "caseStatus in (`$param)"
$ws.Range("E7").Value = "contractLocale in (`$param)"
$ws.Range("F7").Value = "form.setRequired(`$param);"

# Row 9 data values (now quoted, comma separated list of quoted tokens)
$ws.Range("B9").Value = '"BRELETTER1","DLVR1EX1","UWREPL"'
$ws.Range("C9").Value = '"INTERNET"'
$ws.Range("D9").Value = '"SUBMITTED"'
$ws.Range("E9").Value = '"MN"'

# Row 10 data values
$ws.Range("B10").Value = '"PSTDECLINE"'
$ws.Range("C10").Value = '"CAPTIVE","IPAPER","INTERNET"'
$ws.Range("D10").Value = '"AUTODECLINE"'
$ws.Range("E10").Value = '"AZ","CA","CT","GA","IL","MN"'

# Row 11 data values
$ws.Range("B11").Value = '"BRELETTER1","DLVR1EX1","DLVR1APPCOPY","UWREPL"'
$ws.Range("C11").Value = '"CAPTIVE","IPAPER"'
$ws.Range("D11").Value = '"SUBMITTED"'
$ws.Range("E11").Value = '"GA"'

# Column B got a bit wider
$ws.Columns.Item(2).ColumnWidth = 34

# Selection moved to E11
$ws.Range("E11").Select()
